$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab query (B2): drop the Cohort lookup/column (timing issue fix) ---
# The query previously joined an OPTIONAL MATCH (co:cohort) and returned its
# description as a `Cohort` column; that trailing RETURN item (and the comma
# that preceded it) is removed, everything else in the query is untouched.
$orig = $ws.Range("B2").Value2
$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
if ($orig.EndsWith($cohortSuffix)) {
    $ws.Range("B2").Value2 = $orig.Substring(0, $orig.Length - $cohortSuffix.Length)
}

# --- Sheet view: selection/top-left now focuses B2 instead of C4:E4 / B4 ---
$ws.Range("B2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2

# --- Row height: row 2 shrinks now that the Cohort line is gone (rows 3/4 untouched) ---
$ws.Rows(2).RowHeight = 304.5
